# Fruta / hortaliza, semanal
#
# The underlying data rows (2..38) get re-shuffled: for every row r, the
# "observation" columns (Fecha, Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion, Origen, Precio $/Kg, Kg/unidad) are replaced
# by the values that originally lived in another row, per a fixed
# permutation derived from the target workbook. The identifying columns
# (Mercado ID/Mercado/Region/Codreg/Tipo/Producto/Categoria/Variedad/
# Calidad in columns A,B,C,E,F,G,H,I,J,K,L) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New-row -> Source-row mapping (1-based worksheet rows).
$rowMap = @{
    2 = 12;  3 = 15;  4 = 13;  5 = 38;  6 = 10;  7 = 27;  8 = 23;  9 = 11; 10 = 28;
    11 = 33; 12 = 25; 13 = 29; 14 = 2;  15 = 16; 16 = 30; 17 = 5;  18 = 24; 19 = 34;
    20 = 8;  21 = 14; 22 = 17; 23 = 37; 24 = 6;  25 = 7;  26 = 19; 27 = 21; 28 = 31;
    29 = 32; 30 = 20; 31 = 36; 32 = 9;  33 = 3;  34 = 18; 35 = 4;  36 = 26; 37 = 22;
    38 = 35
}

# Columns that move together as part of the shuffled observation.
$cols = @("D", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot every relevant cell's current value before writing anything,
# so that reads are never affected by in-progress writes.
$snapshot = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 38; $r++) {
        $snapshot["$col$r"] = $ws.Range("$col$r").Value2
    }
}

foreach ($col in $cols) {
    foreach ($newRow in $rowMap.Keys) {
        $srcRow = $rowMap[$newRow]
        $ws.Range("$col$newRow").Value = $snapshot["$col$srcRow"]
    }
}
